# "Starting with numpy and OOP implementation"
# Re-sort the two lookup tables by their first column and refresh the
# active selections, matching the manual re-sort the author performed
# in Excel before committing.

$wb = $excel.ActiveWorkbook

# --- "Course Schedule No Duplicates" sheet: Table167, sort by column A (DISCIPLINA code) ---
$wsDup = $wb.Worksheets.Item(2)
$wsDup.Activate() | Out-Null

$tblDup = $wsDup.ListObjects.Item(1)
$sortDup = $tblDup.Sort
$sortDup.SortFields.Clear() | Out-Null
$sortDup.SortFields.Add($tblDup.ListColumns.Item(1).Range) | Out-Null
$sortDup.Apply() | Out-Null

# --- "Courses Of Each Professor" sheet: Table2, sort by column B (DISCIPLINA code) ---
$wsProf = $wb.Worksheets.Item(3)
$wsProf.Activate() | Out-Null

$tblProf = $wsProf.ListObjects.Item(1)
$sortProf = $tblProf.Sort
$sortProf.SortFields.Clear() | Out-Null
$sortProf.SortFields.Add($tblProf.ListColumns.Item(1).Range) | Out-Null
$sortProf.Apply() | Out-Null

# Leave the final selection as it was left by the author: A2 on the
# "Course Schedule No Duplicates" sheet (still the active tab) and B8
# on "Courses Of Each Professor".
$wsProf.Range("B8").Select() | Out-Null

$wsDup.Activate() | Out-Null
$wsDup.Range("A2").Select() | Out-Null
